$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows just before the "EndOfData" marker row (row 23),
# pushing it down to row 25.
$ws.Rows("23:24").Insert()

# Copy the formatting of the last data row ("Category 5", row 21) into the
# two freshly inserted rows so the new entries look like the rest of the
# table (same fonts/fills/borders per column).
$ws.Range("A21:H21").Copy()
$ws.Range("A22:H23").PasteSpecial(-4122)  # xlPasteFormats

# Row 22: "Data Scheme" test pattern
$ws.Range("C22").Value = "Data Scheme"
$ws.Range("G22").Value = "data:datascheme"
$ws.Range("H22").Value = 1

# Row 23: "JavaScript Scheme" test pattern
$ws.Range("C23").Value = "JavaScript Scheme"
$ws.Range("G23").Value = "javascript:alert(123);"
$ws.Range("H23").Value = 1

# Hyperlinks for the two new "path" cells
$ws.Hyperlinks.Add($ws.Range("G22"), "data:datascheme", "", "", "data:datascheme")
$ws.Hyperlinks.Add($ws.Range("G23"), "javascript:alert(123);", "", "", "javascript:alert(123")

# Restore the selection state that Excel would persist on save
$ws.Range("H33").Select()
